# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price/volume columns are plain text cells (e.g. "65.305.40" isn't a valid
# number), so numeric-looking values are written with a leading apostrophe
# to keep Excel from auto-converting them to numbers, exactly like a human
# typing '568.97 into a General-formatted cell would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.305.40"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "2.937.35"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'568.97"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").Value = "'159.20"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "2.933.20"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").Value = "'6.74"
$ws.Range("E10").Value = "  -3.53%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").Value = "'34.56"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D16").Value = "65.334.83"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "3.424.67"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").Value = "'7.05"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "2.936.18"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").Value = "'15.51"
$ws.Range("E20").Value = "  +11.51%  "
$ws.Range("D21").Value = "'445.13"
$ws.Range("E21").Value = "  -4.05%  "
$ws.Range("D22").Value = "'0.697"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").Value = "'7.29"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "'82.39"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("E27").Value = "  -5.77%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "'2.39"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "0.0₃0999"
$ws.Range("E32").Value = "  -5.77%  "
$ws.Range("D33").Value = "'27.18"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'0.973"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").Value = "'5.75"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'44.28"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  -8.66%  "
$ws.Range("D41").Value = "'0.303"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "'2.84"
$ws.Range("E43").Value = "  -7.55%  "
$ws.Range("D44").Value = "'8.51"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").Value = "'382.37"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "2.694.51"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").Value = "'133.33"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +4.72%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'23.40"
$ws.Range("E51").Value = "  -0.71%  "
